$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The uploaded workbook re-sorts the hourly "접수시간대" (reception-hour)
# frequency table into descending order by count and refreshes the counts.
# Row 1 (header) is left untouched. Rows 2-25 get the new hour labels (A)
# and new counts (B). The hour labels are text (e.g. "08", "09") so the
# cells are pre-formatted as Text to keep Excel from dropping the leading
# zero / turning them into numbers.
$labels = @("08","11","09","10","17","15","16","13","24","14","07","18","19","20","21","22","23","05","03","01","06","04","02","12")
$values = @(72,59,58,58,42,41,40,40,33,28,27,22,17,13,12,8,7,5,4,3,2,1,1,1)

$ws.Range("A2:A25").NumberFormat = "@"

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
